# Update cryptocurrency price/volume data per latest GitHub Actions scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''27.568.75'
$ws.Range('E2').Value = '  +2.52%  '

# Row 3
$ws.Range('D3').Value = '''1.852.83'
$ws.Range('E3').Value = '  +2.18%  '

# Row 4
$ws.Range('D4').Value = '''1.031'
$ws.Range('E4').Value = '  +2.74%  '

# Row 5
$ws.Range('D5').Value = '''321.43'
$ws.Range('E5').Value = '  +3.20%  '

# Row 6
$ws.Range('D6').Value = '''1.028'
$ws.Range('E6').Value = '  +2.45%  '

# Row 7
$ws.Range('D7').Value = '''0.4387'
$ws.Range('E7').Value = '  +2.28%  '

# Row 8
$ws.Range('D8').Value = '''0.3776'
$ws.Range('E8').Value = '  +2.34%  '

# Row 9
$ws.Range('D9').Value = '''0.07413'
$ws.Range('E9').Value = '  +2.40%  '

# Row 10
$ws.Range('D10').Value = '''0.8764'
$ws.Range('E10').Value = '  +1.69%  '

# Row 11
$ws.Range('D11').Value = '''21.53'
$ws.Range('E11').Value = '  +1.58%  '

# Row 12
$ws.Range('D12').Value = '''1.854.19'
$ws.Range('E12').Value = '  -8.00%  '

# Row 13
$ws.Range('D13').Value = '''5.529'
$ws.Range('E13').Value = '  +2.47%  '

# Row 14
$ws.Range('D14').Value = '''6.702'
$ws.Range('E14').Value = '  +0.96%  '

# Row 15
$ws.Range('D15').Value = '''0.07210'
$ws.Range('E15').Value = '  +4.57%  '

# Row 16
$ws.Range('D16').Value = '''82.94'
$ws.Range('E16').Value = '  +2.79%  '

# Row 17
$ws.Range('E17').Value = '  +2.45%  '

# Row 18
$ws.Range('D18').Value = '''0.000009028'
$ws.Range('E18').Value = '  +1.06%  '

# Row 19
$ws.Range('D19').Value = '''1.028'
$ws.Range('E19').Value = '  +2.46%  '

# Row 20
$ws.Range('D20').Value = '''15.42'
$ws.Range('E20').Value = '  +1.63%  '

# Row 21
$ws.Range('D21').Value = '''27.573.39'
$ws.Range('E21').Value = '  +2.35%  '

# Row 22
$ws.Range('D22').Value = '''5.258'
$ws.Range('E22').Value = '  +1.59%  '

# Row 23
$ws.Range('D23').Value = '''11.35'
$ws.Range('E23').Value = '  +2.72%  '

# Row 24
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').Value = '''157.84'
$ws.Range('E24').Value = '  +2.78%  '

# Row 25
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '''1.914'
$ws.Range('E25').Value = '  +1.72%  '

# Row 26
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '''18.77'
$ws.Range('E26').Value = '  +2.80%  '

# Row 27
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '''1.973'
$ws.Range('E27').Value = '  +5.69%  '

# Row 28
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '''5.265'
$ws.Range('E28').Value = '  +0.77%  '

# Row 29
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = '''117.05'
$ws.Range('E29').Value = '  +1.87%  '

# Row 30
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value = '''0.09055'
$ws.Range('E30').Value = '  +1.30%  '

# Row 31
$ws.Range('B31').Value = 'ARBITRUM'
$ws.Range('C31').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D31').Value = '''1.198'
$ws.Range('E31').Value = '  +2.79%  '

# Row 32
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '''0.7627'
$ws.Range('E32').Value = '  +2.56%  '

# Row 33
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '''4.521'
$ws.Range('E33').Value = '  +2.25%  '

# Row 34
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '''2.873'
$ws.Range('E34').Value = '  +2.75%  '

# Row 35
$ws.Range('B35').Value = 'Frax'
$ws.Range('C35').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D35').Value = '''1.029'
$ws.Range('E35').Value = '  +2.15%  '

# Row 36
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').Value = '''1.149'
$ws.Range('E36').Value = '  +2.90%  '

# Row 37
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.01978'
$ws.Range('E37').Value = '  +3.02%  '

# Row 38
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '''0.05303'
$ws.Range('E38').Value = '  +1.77%  '

# Row 39
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').Value = '''0.5155'
$ws.Range('E39').Value = '  +1.51%  '

# Row 40
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '''2.801'
$ws.Range('E40').Value = '  +3.17%  '

# Row 41
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = '''0.1676'
$ws.Range('E41').Value = '  +2.21%  '

# Row 42
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '''6.739'
$ws.Range('E42').Value = '  +4.72%  '

# Row 43
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '''8.485'
$ws.Range('E43').Value = '  +2.85%  '

# Row 44
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '''108.75'
$ws.Range('E44').Value = '  +1.87%  '

# Row 45
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '''10.52'
$ws.Range('E45').Value = '  +1.61%  '

# Row 46
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').Value = '''1.710'
$ws.Range('E46').Value = '  +3.41%  '

# Row 47
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '''0.4652'
$ws.Range('E47').Value = '  +2.18%  '

# Row 48
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.06394'
$ws.Range('E48').Value = '  +1.74%  '

# Row 49
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '''1.858'
$ws.Range('E49').Value = '  +2.92%  '

# Row 50
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = '''39.13'
$ws.Range('E50').Value = '  +4.04%  '

# Row 51
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '''63.97'
$ws.Range('E51').Value = '  +1.26%  '
